$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header "Label" in column H, matching the bold/bordered/centered
# style already used by the other header cells (B1:G1)
$ws.Range("H1").Value = "Label"
$ws.Range("H1").Font.Bold = $true
$ws.Range("H1").HorizontalAlignment = -4108
$ws.Range("H1").VerticalAlignment = -4160
$ws.Range("H1").Borders.LineStyle = 1
$ws.Range("H1").Borders.Weight = 2

# Updated D/E refit values and new H (Label) column values
$ws.Range("D2").Value = 0.443252525333779
$ws.Range("E2").Value = 0.443252525333779
$ws.Range("H2").Value = 0

$ws.Range("D3").Value = 0.2949300898600171
$ws.Range("E3").Value = 0.2949300898600171
$ws.Range("H3").Value = 0

$ws.Range("D4").Value = 0.5125453428804139
$ws.Range("E4").Value = 0.5125453428804139
$ws.Range("H4").Value = 0

$ws.Range("D5").Value = 0.2161368231992738
$ws.Range("E5").Value = 0.2161368231992738
$ws.Range("H5").Value = 0

$ws.Range("D6").Value = 0.5713385061223186
$ws.Range("E6").Value = 0.5713385061223186
$ws.Range("H6").Value = 0

$ws.Range("D7").Value = 0.5083698670642057
$ws.Range("E7").Value = 0.4916301329357943
$ws.Range("H7").Value = 1

$ws.Range("D8").Value = 0.5011080771627284
$ws.Range("E8").Value = 0.4988919228372716
$ws.Range("H8").Value = 1

$ws.Range("D9").Value = 0.5024697877591983
$ws.Range("E9").Value = 0.4975302122408017
$ws.Range("H9").Value = 1

$ws.Range("D10").Value = 0.4648805946301293
$ws.Range("E10").Value = 0.5351194053698707
$ws.Range("H10").Value = 1

$ws.Range("D11").Value = 0.1704782088111302
$ws.Range("E11").Value = 0.8295217911888697
$ws.Range("H11").Value = 1

$ws.Range("H12").Value = 0
$ws.Range("H13").Value = 0
$ws.Range("H14").Value = 0
$ws.Range("H15").Value = 0
$ws.Range("H16").Value = 0
$ws.Range("H17").Value = 1
$ws.Range("H18").Value = 1
$ws.Range("H19").Value = 1
$ws.Range("H20").Value = 1
$ws.Range("H21").Value = 1
